$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)
$hf = $sec.Headers.Item(2)
$shp = $hf.Shapes.Item(3)
$tf = $shp.TextFrame
$tr = $tf.TextRange
"Start/End: " + $tr.Start + "/" + $tr.End
"StoryType: " + $tr.StoryType
# try moving within the range
try {
  $tr2 = $tr.Duplicate
  "duplicate ok, text=[" + $tr2.Text + "]"
} catch {
  "duplicate err: " + $_
}
# Try Find with minimal args
try {
    $f = $tr.Find
    $f.ClearFormatting()
    $f.Text = "Lâm"
    $res = $f.Execute()
    "Find.Execute (prop-based) result: " + $res
} catch {
   "err2: " + $_
}
